# Adapt tests to control version
# Add a new "version" column to the "settings" sheet of the XLSForm
# workbook, mirroring form_title / form_id, and set its value to 1 for
# the single settings row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Header: C1 = "version" (same row as form_title / form_id headers)
$ws.Range("C1").Value = "version"

# Data: C2 = 1 (the control version number for this settings row)
$ws.Range("C2").Value = 1

# Leave the new cell selected, matching the cursor resting just past the
# newly entered data after typing it in.
[void]$ws.Range("C3").Select()
